# Update countries & provincias Spain
# Applies the 19-May-2020 00:35 -> 01:05 refresh of the COVID-19 country
# table: updated case counters for several countries plus the consequent
# re-ranking (swaps) of a handful of neighbouring rows, and the "last
# updated" timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 01:05"

# --- Row 4: Estados Unidos (rank unchanged, counters refreshed) ---
$ws.Cells.Item(4, 2).Value = 1548357
$ws.Cells.Item(4, 3).Value = 20693
$ws.Cells.Item(4, 4).Value = 353781
$ws.Cells.Item(4, 5).Value = 1102717
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 881
$ws.Cells.Item(4, 8).Value = 91859

# --- Rows 7-8: Brasil overtakes Reino Unido ---
$ws.Cells.Item(7, 1).Value = "Brasil"
$ws.Cells.Item(7, 2).Value = 254220
$ws.Cells.Item(7, 3).Value = 13140
$ws.Cells.Item(7, 4).Value = 100459
$ws.Cells.Item(7, 5).Value = 136969
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 674
$ws.Cells.Item(7, 8).Value = 16792

$ws.Cells.Item(8, 1).Value = "Reino Unido"
$ws.Cells.Item(8, 2).Value = 246406
$ws.Cells.Item(8, 3).Value = 2711
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 160
$ws.Cells.Item(8, 8).Value = 34796

# --- Row 26: Suiza (active/recovered refreshed only) ---
$ws.Cells.Item(26, 4).Value = 27600
$ws.Cells.Item(26, 5).Value = 1111

# --- Row 51: Chequia (rank unchanged, counters refreshed) ---
$ws.Cells.Item(51, 2).Value = 8586
$ws.Cells.Item(51, 3).Value = 111
$ws.Cells.Item(51, 4).Value = 5641
$ws.Cells.Item(51, 5).Value = 2648
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 297

# --- Rows 52-53: Argentina overtakes Noruega ---
$ws.Cells.Item(52, 1).Value = "Argentina"
$ws.Cells.Item(52, 2).Value = 8371
$ws.Cells.Item(52, 3).Value = 303
$ws.Cells.Item(52, 4).Value = 2625
$ws.Cells.Item(52, 5).Value = 5364
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 9
$ws.Cells.Item(52, 8).Value = 382

$ws.Cells.Item(53, 1).Value = "Noruega"
$ws.Cells.Item(53, 2).Value = 8257
$ws.Cells.Item(53, 3).Value = 8
$ws.Cells.Item(53, 4).Value = 32
$ws.Cells.Item(53, 5).Value = 7992
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 1
$ws.Cells.Item(53, 8).Value = 233

# --- Row 111: Niger (counters refreshed) ---
$ws.Cells.Item(111, 2).Value = 909
$ws.Cells.Item(111, 3).Value = 5
$ws.Cells.Item(111, 4).Value = 714
$ws.Cells.Item(111, 5).Value = 140
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 55

# --- Row 112: Uruguay (B-E refreshed only) ---
$ws.Cells.Item(112, 2).Value = 887
$ws.Cells.Item(112, 3).Value = 153
$ws.Cells.Item(112, 4).Value = 569
$ws.Cells.Item(112, 5).Value = 298

# --- Rows 126-130: Haiti jumps ahead of Jamaica / Chad / Sierra Leona / Tanzania ---
$ws.Cells.Item(126, 1).Value = "Haiti"
$ws.Cells.Item(126, 2).Value = 533
$ws.Cells.Item(126, 3).Value = 77
$ws.Cells.Item(126, 4).Value = 21
$ws.Cells.Item(126, 5).Value = 491
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 21

$ws.Cells.Item(127, 1).Value = "Jamaica"
$ws.Cells.Item(127, 2).Value = 520
$ws.Cells.Item(127, 3).Value = 3
$ws.Cells.Item(127, 4).Value = 127
$ws.Cells.Item(127, 5).Value = 384
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 9

$ws.Cells.Item(128, 1).Value = "Republica del Chad"
$ws.Cells.Item(128, 2).Value = 519
$ws.Cells.Item(128, 3).Value = 16
$ws.Cells.Item(128, 4).Value = 117
$ws.Cells.Item(128, 5).Value = 349
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 53

$ws.Cells.Item(129, 1).Value = "Sierra Leona"
$ws.Cells.Item(129, 2).Value = 519
$ws.Cells.Item(129, 3).Value = 14
$ws.Cells.Item(129, 4).Value = 148
$ws.Cells.Item(129, 5).Value = 338
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 1
$ws.Cells.Item(129, 8).Value = 33

$ws.Cells.Item(130, 1).Value = "Tanzania"
$ws.Cells.Item(130, 2).Value = 509
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 183
$ws.Cells.Item(130, 5).Value = 305
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 21

# --- Rows 162-163: Bermudas overtakes Guyana ---
$ws.Cells.Item(162, 1).Value = "Bermudas"
$ws.Cells.Item(162, 2).Value = 125
$ws.Cells.Item(162, 3).Value = 2
$ws.Cells.Item(162, 4).Value = 77
$ws.Cells.Item(162, 5).Value = 39
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 9

$ws.Cells.Item(163, 1).Value = "Guyana"
$ws.Cells.Item(163, 2).Value = 124
$ws.Cells.Item(163, 3).Value = 7
$ws.Cells.Item(163, 4).Value = 44
$ws.Cells.Item(163, 5).Value = 70
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 10
